$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - "Förändrad" date bump
$ws.Range("C2").Value = 45180

# Row 2 - updated species counts
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 5
$ws.Range("K2").Value = 2
$ws.Range("O2").Value = 7
$ws.Range("P2").Value = 2
$ws.Range("Q2").Value = 17

# Row 2 - updated species name list
$r2Text = "Goliatmusseron`r`nLakritsmusseron`r`nBlå taggsvamp`r`nSpillkråka`r`nSvart taggsvamp`r`nSvartvit taggsvamp`r`nTallticka`r`nBjörksplintborre`r`nBlåmossa`r`nBronshjon`r`nDropptaggsvamp`r`nGranbarkgnagare`r`nGrovticka`r`nMindre märgborre`r`nSkarp dropptaggsvamp`r`nFläcknycklar`r`nRevlummer"
$ws.Range("R2").Value = $r2Text

# "Förändrad" date bump for remaining rows
$ws.Range("C3").Value = 45180
$ws.Range("C4").Value = 45180
$ws.Range("C5").Value = 45180
$ws.Range("C6").Value = 45180
$ws.Range("C7").Value = 45180
$ws.Range("C8").Value = 45180
$ws.Range("C9").Value = 45180
$ws.Range("C10").Value = 45180
$ws.Range("C11").Value = 45180
